# Updated symbol list refresh (matches commit: 'Updated symbol list on Tue Feb 14 15:14:31 UTC 2023 with GitHub Actions').
# Rows 6-17 rotate down by one coin (new entry GateToken enters at the top of that block,
# BitKan's block drops SpecialPowerGold... i.e. the ranking list shifted), and Price/Volume(1h)/Hora
# refresh for (almost) every row. Hora (col G) moves 14 -> 15 for all data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry lists only the columns that actually change for that row.
# D/E/G values are prefixed with a leading apostrophe so Excel stores them as literal
# text (matching the workbook's existing text-cell convention) instead of auto-coercing
# numeric-looking strings ('297.91', '2.60%', '15') into Number/Percentage cells.
$updates = @(
    @{Row=2; D='297.91'; E='2.60%'; G='15'},
    @{Row=3; D='41.46'; E='4.29%'; G='15'},
    @{Row=4; D='5.048'; E='0.41%'; G='15'},
    @{Row=5; D='0.07482'; E='1.70%'; G='15'},
    @{Row=6; B='GateToken'; C='https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'; D='4.342'; E='1.27%'; G='15'},
    @{Row=7; B='FTXToken'; C='https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'; D='1.569'; E='1.18%'; G='15'},
    @{Row=8; B='MXToken'; C='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D='0.9379'; E='2.80%'; G='15'},
    @{Row=9; B='BTSEToken'; C='https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'; D='2.401'; E='0.17%'; G='15'},
    @{Row=10; B='LiechtensteinCryptoassetsExchange'; C='https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'; D='0.1204'; E='0.98%'; G='15'},
    @{Row=11; B='WazirX'; C='https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'; D='0.1818'; E='3.52%'; G='15'},
    @{Row=12; B='MandalaExchangeToken'; C='https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'; D='0.08800'; E='0.80%'; G='15'},
    @{Row=13; B='BitrueCoin'; C='https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'; D='0.04239'; E='1.85%'; G='15'},
    @{Row=14; B='BitMartToken'; C='https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'; D='0.1048'; E='-0.15%'; G='15'},
    @{Row=15; B='BitForexToken'; C='https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'; D='0.001265'; E='-1.05%'; G='15'},
    @{Row=16; B='TigerCash'; C='https://coinranking.com/coin/6hIn06L2+tigercash-tch'; D='0.005857'; E='1.05%'; G='15'},
    @{Row=17; B='LEO'; C='https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'; D='3.352'; E='-1.27%'; G='15'},
    @{Row=18; D='0.3306'; E='0.53%'; G='15'},
    @{Row=19; D='7.865'; E='4.19%'; G='15'},
    @{Row=20; D='0.1373'; E='1.52%'; G='15'},
    @{Row=21; D='0.3276'; E='13.57%'; G='15'},
    @{Row=22; D='0.03964'; E='3.20%'; G='15'},
    @{Row=23; D='0.001258'; E='-0.97%'; G='15'},
    @{Row=24; D='0.003856'; E='-0.84%'; G='15'},
    @{Row=25; D='0.0001221'; E='-4.76%'; G='15'},
    @{Row=26; D='0.0003695'; E='-0.90%'; G='15'},
    @{Row=27; G='15'},
    @{Row=28; G='15'},
    @{Row=29; G='15'},
    @{Row=30; G='15'},
    @{Row=31; G='15'},
    @{Row=32; G='15'},
    @{Row=33; G='15'},
    @{Row=34; G='15'},
    @{Row=35; G='15'},
    @{Row=36; G='15'},
    @{Row=37; G='15'},
    @{Row=38; D='0.02407'; E='3.02%'; G='15'},
    @{Row=39; D='0.05155'; E='2.50%'; G='15'},
    @{Row=40; D='0.005887'; E='15.21%'; G='15'},
    @{Row=41; D='0.007714'; E='0.22%'; G='15'},
    @{Row=42; D='0.1322'; E='3.82%'; G='15'},
    @{Row=43; D='0.007323'; E='-0.67%'; G='15'},
    @{Row=44; D='0.007125'; E='2.12%'; G='15'},
    @{Row=45; D='0.2961'; E='-6.02%'; G='15'},
    @{Row=46; D='0.00006183'; E='-5.35%'; G='15'},
    @{Row=47; D='0.00000000744'; E='-0.90%'; G='15'},
    @{Row=48; D='0.04626'; E='-81.63%'; G='15'},
    @{Row=49; D='0.004169'; E='-0.89%'; G='15'},
    @{Row=50; D='0.00002084'; E='-0.90%'; G='15'},
    @{Row=51; D='0.0001985'; E='-0.90%'; G='15'},
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey('B')) { $ws.Cells.Item($r, 2).Value = $u.B }
    if ($u.ContainsKey('C')) { $ws.Cells.Item($r, 3).Value = $u.C }
    if ($u.ContainsKey('D')) { $ws.Cells.Item($r, 4).Value = "'" + $u.D }
    if ($u.ContainsKey('E')) { $ws.Cells.Item($r, 5).Value = "'" + $u.E }
    if ($u.ContainsKey('G')) { $ws.Cells.Item($r, 7).Value = "'" + $u.G }
}

